$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Ativação date bump: 01/01/2018 -> 01/01/2022
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2018", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2022", 2)

# ---------------------------------------------------------------------
# Helper: find the 1-based index of the paragraph whose text starts
# with $prefix (searched fresh each time so indices stay correct even
# after earlier insertions shift everything below them).
# ---------------------------------------------------------------------
function Find-ParagraphIndex($doc, $prefix) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# Helper: insert a new italic paragraph right after the paragraph that
# starts with $afterPrefix, with the given $text.
function Insert-ItalicParagraphAfter($doc, $afterPrefix, $text) {
    $idx = Find-ParagraphIndex $doc $afterPrefix
    $src = $doc.Paragraphs.Item($idx)
    $src.Range.InsertParagraphAfter()
    $newPara = $doc.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $text
    $newRange = $doc.Range($newPara.Range.Start, $newPara.Range.Start + $text.Length)
    $newRange.Italic = 1
}

# ---------------------------------------------------------------------
# 2) Objetivos: add italic English translation paragraph
# ---------------------------------------------------------------------
$objEn = "Provide students with knowledge about environmental management in companies, environmental policies, environmental management systems (EMS) and ISO 14000 series standards, enabling them to participate in the planning and implementation of an EMS in a company."
Insert-ItalicParagraphAfter $d "Proporcionar aos alunos conhecimentos sobre a gestão ambiental" $objEn

# ---------------------------------------------------------------------
# 3) Programa resumido: add italic English translation paragraph
# ---------------------------------------------------------------------
$sumEn = "Environmental Management Systems; Iso 14000; Environmental Audit."
Insert-ItalicParagraphAfter $d "Sistemas de Gestão Ambiental; Iso 14000; Auditoria Ambiental." $sumEn

# ---------------------------------------------------------------------
# 4) Programa: rewrite Portuguese body text, then add italic English
#    translation paragraph after it.
# ---------------------------------------------------------------------
$progPt = "Evolução das práticas de gestão ambiental empresarial;- Economia circular, conceitos e aplicações;- Responsabilidade social corporativa: conceito e programa;- Implantação do sistema de gerenciamento ambiental (SGA): conceitos e modelos;- Produção mais limpa;- Ferramentas de gestão focadas no produto;- Análise e otimização do ciclo de vida do produto;- Ecoinovação e Ecodesign;- Rotulagem ambiental;- Inovação e sustentabilidade;- Normas ISO 14001 (série ISO 14000), requisitos e orientações para uso e Certificações ambientais."
$progEn = "Evolution of corporate environmental management practices;- Circular economy, concepts and applications;- Corporate social responsibility: concept and program;- Implementation of the environmental management system (SGA): concepts and models;- Cleaner production;- Management tools focused on the product;- Analysis and optimization of the product life cycle;- Eco-innovation and Ecodesign;- Environmental labeling;- Innovation and sustainability;- ISO 14001 standards (ISO 14000 series), requirements and guidelines for use and Environmental Certifications."

$progIdx = Find-ParagraphIndex $d "Evolução das práticas de gestão ambiental empresarial"
$progPara = $d.Paragraphs.Item($progIdx)
$progPara.Range.Text = $progPt

Insert-ItalicParagraphAfter $d "Evolução das práticas de gestão ambiental empresarial;- Economia circular" $progEn

Write-Output "done"
